$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as literal text (to preserve the dotted
# thousands grouping / trailing zeros used by the source site), so for any
# new value that Excel would otherwise auto-convert to a real number we pin
# the cell format to Text ("@") first, keeping the original exact string.

$ws.Range("D2").Value2 = "30.440.99"
$ws.Range("E2").Value2 = "  -1.12%  "

$ws.Range("E3").Value2 = "  +1.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.9999"
$ws.Range("E4").Value2 = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "243.88"
$ws.Range("E5").Value2 = "  +2.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.000"
$ws.Range("E6").Value2 = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4700"
$ws.Range("E7").Value2 = "  -1.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2865"
$ws.Range("E8").Value2 = "  -0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06811"
$ws.Range("E9").Value2 = "  +3.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "109.40"
$ws.Range("E10").Value2 = "  +11.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "18.43"
$ws.Range("E11").Value2 = "  -1.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.07741"
$ws.Range("E12").Value2 = "  +1.90%  "

$ws.Range("D13").Value2 = "1.889.19"
$ws.Range("E13").Value2 = "  +0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.291"
$ws.Range("E14").Value2 = "  +3.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.6584"
$ws.Range("E15").Value2 = "  +0.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "293.97"
$ws.Range("E16").Value2 = "  -4.00%  "

$ws.Range("D17").Value2 = "30.424.21"
$ws.Range("E17").Value2 = "  -1.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.000007622"
$ws.Range("E18").Value2 = "  +0.65%  "

$ws.Range("E19").Value2 = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "12.94"
$ws.Range("E20").Value2 = "  -1.58%  "

$ws.Range("D21").Value2 = "2.132.83"
$ws.Range("E21").Value2 = "  +0.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "1.001"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.248"
$ws.Range("E23").Value2 = "  +2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "6.206"
$ws.Range("E24").Value2 = "  +0.90%  "

$ws.Range("E25").Value2 = "  +7.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "9.337"
$ws.Range("E26").Value2 = "  +0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "168.69"
$ws.Range("E27").Value2 = "  +1.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.083"
$ws.Range("E28").Value2 = "  +7.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.1070"
$ws.Range("E29").Value2 = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.366"
$ws.Range("E30").Value2 = "  +0.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.175"
$ws.Range("E31").Value2 = "  +0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.981"
$ws.Range("E32").Value2 = "  +0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.05045"
$ws.Range("E33").Value2 = "  +0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.7387"
$ws.Range("E34").Value2 = "  +1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.155"
$ws.Range("E35").Value2 = "  -1.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.02065"
$ws.Range("E36").Value2 = "  +6.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.678"
$ws.Range("E38").Value2 = "  -0.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.061"
$ws.Range("E39").Value2 = "  -0.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "109.27"
$ws.Range("E40").Value2 = "  +1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.8726"
$ws.Range("E41").Value2 = "  -3.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "5.829"
$ws.Range("E42").Value2 = "  +3.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.4267"
$ws.Range("E43").Value2 = "  +1.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.9999"
$ws.Range("E44").Value2 = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "51.44"
$ws.Range("E45").Value2 = "  +20.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "67.38"
$ws.Range("E46").Value2 = "  +2.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "7.199"
$ws.Range("E47").Value2 = "  -1.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "9.246"
$ws.Range("E48").Value2 = "  +2.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.1217"
$ws.Range("E49").Value2 = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "34.97"
$ws.Range("E50").Value2 = "  +0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.2446"
$ws.Range("E51").Value2 = "  +11.32%  "
